$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 1144, shifting all rows
# from 1144..1180 down to 1147..1183 (dates/prices/volumes for the newest
# reporting week land at the top of the historical block).
$ws.Rows("1144:1146").Insert()

# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg,
# F Categoría ID, G Categoría, H Variedad, I Calidad, J Volumen,
# K Precio mínimo, L Precio máximo, M Precio promedio ponderado,
# N Unidad de comercialización, O Origen, P Precio $/Kg,
# Q Kg o Unidades, R Clasificación
$newRows = @(
    @(1144, "Lluteño", "Primera", 40, 44000, 45000, 44500, "`$/saco 50 unidades", 890, 50),
    @(1145, "Lluteño", "Segunda", 50, 39000, 40000, 39500, "`$/saco 75 unidades", 527, 75),
    @(1146, "Lluteño", "Tercera", 50, 34000, 35000, 34500, "`$/saco 100 unidades", 345, 100)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = 1
    $ws.Cells.Item($rowNum, 2).Value = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($rowNum, 3).Value = "Arica y Parinacota"
    $ws.Cells.Item($rowNum, 4).Value = 45239
    $ws.Cells.Item($rowNum, 5).Value = 15
    $ws.Cells.Item($rowNum, 6).Value = 100112024
    $ws.Cells.Item($rowNum, 7).Value = "Choclo"
    $ws.Cells.Item($rowNum, 8).Value = $r[1]
    $ws.Cells.Item($rowNum, 9).Value = $r[2]
    $ws.Cells.Item($rowNum, 10).Value = $r[3]
    $ws.Cells.Item($rowNum, 11).Value = $r[4]
    $ws.Cells.Item($rowNum, 12).Value = $r[5]
    $ws.Cells.Item($rowNum, 13).Value = $r[6]
    $ws.Cells.Item($rowNum, 14).Value = $r[7]
    $ws.Cells.Item($rowNum, 15).Value = "Región de Arica y Parinacota"
    $ws.Cells.Item($rowNum, 16).Value = $r[8]
    $ws.Cells.Item($rowNum, 17).Value = $r[9]
    $ws.Cells.Item($rowNum, 18).Value = "Hortaliza"
}
